$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.705.89"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.051.71"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.59"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.20"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.043.59"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.21"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000223"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "3.547.94"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "62.748.50"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "3.053.91"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.66"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.47"
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.20"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.94"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.74"
$ws.Range("E33").Value = "  +6.42%  "
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.10"
$ws.Range("E35").Value = "  -6.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "459.04"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").Value = "3.193.69"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.18"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.72"
$ws.Range("E44").Value = "  +7.48%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.247"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.99"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.72"
$ws.Range("E49").Value = "  -5.05%  "
$ws.Range("D50").Value = "0.0₃0496"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("E51").Value = "  +3.26%  "
